# "final for 2E python section"
# The price/points stat tables were produced with duplicated pandas column
# names ("Chile price price" / "_1" / "Chile points points"). This cleans
# that up: the empty spacer column C ("_1") is removed entirely (shifting
# the points column from D into C), and the two header labels are
# collapsed into their corrected single-token form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stray spacer column (old C, header "_1") - this shifts the
# "Chile points points" column from D to C and tightens the used range
# from A1:D9 down to A1:C9.
$ws.Columns("C").Delete()

# Fix up the (now B1/C1) header labels.
$ws.Range("B1").Value = "Chile_priceprice"
$ws.Range("C1").Value = "Chile_pointspoints"
